$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its values as text, matching the
# original inline-string cell type. Without this, values like "1.00" or
# "485.27" would be auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.046.73"
$ws.Range("E2").Value = "  +1.72%  "

$ws.Range("D3").Value = "3.910.43"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "485.27"
$ws.Range("E5").Value = "  +3.44%  "

$ws.Range("D6").Value = "145.75"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("D8").Value = "0.997"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").Value = "0.728"
$ws.Range("E9").Value = "  -1.86%  "

$ws.Range("E10").Value = "  +1.33%  "

$ws.Range("D11").Value = "0.0000355"
$ws.Range("E11").Value = "  +4.91%  "

$ws.Range("D12").Value = "42.51"
$ws.Range("E12").Value = "  -1.59%  "

$ws.Range("D13").Value = "10.64"

$ws.Range("D14").Value = "4.529.15"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").Value = "14.73"
$ws.Range("E15").Value = "  -1.79%  "

$ws.Range("D16").Value = "3.899.01"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("D18").Value = "19.83"
$ws.Range("E18").Value = "  -0.95%  "

$ws.Range("E19").Value = "  -2.45%  "

$ws.Range("D20").Value = "68.131.86"
$ws.Range("E20").Value = "  +1.36%  "

$ws.Range("D21").Value = "448.86"
$ws.Range("E21").Value = "  +3.80%  "

$ws.Range("D22").Value = "14.69"
$ws.Range("E22").Value = "  -0.36%  "

$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("D24").Value = "88.97"
$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("D25").Value = "11.71"
$ws.Range("E25").Value = "  +15.52%  "

$ws.Range("D26").Value = "11.05"
$ws.Range("E26").Value = "  +14.02%  "

$ws.Range("D27").Value = "3.61"
$ws.Range("E27").Value = "  +1.99%  "

$ws.Range("D28").Value = "38.73"
$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("D29").Value = "5.84"
$ws.Range("E29").Value = "  +3.31%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.131"
$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "689.36"
$ws.Range("E31").Value = "  -6.26%  "

$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "13.36"
$ws.Range("E32").Value = "  -2.49%  "

$ws.Range("D33").Value = "2.86"
$ws.Range("E33").Value = "  +3.01%  "

$ws.Range("D34").Value = "0.0₃0935"
$ws.Range("E34").Value = "  +25.52%  "

$ws.Range("D35").Value = "41.64"
$ws.Range("E35").Value = "  -5.33%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "5.80"
$ws.Range("E36").Value = "  +8.26%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "58.74"
$ws.Range("E37").Value = "  +1.09%  "

$ws.Range("D38").Value = "0.150"
$ws.Range("E38").Value = "  -5.34%  "

$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("D40").Value = "0.0477"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").Value = "2.88"
$ws.Range("E41").Value = "  +16.01%  "

$ws.Range("D42").Value = "3.04"
$ws.Range("E42").Value = "  -7.17%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "0.359"
$ws.Range("E43").Value = "  +6.76%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "3.02"
$ws.Range("E44").Value = "  +7.62%  "

$ws.Range("E45").Value = "  +0.49%  "

$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("D48").Value = "2.13"
$ws.Range("E48").Value = "  -2.95%  "

$ws.Range("D49").Value = "146.52"
$ws.Range("E49").Value = "  +2.30%  "

$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  -1.45%  "

$ws.Range("D51").Value = "0.0₆0333"
$ws.Range("E51").Value = "  +44.65%  "
